# Update "想去人数" (column F) counts for a batch of events.
# Sheet 1 = 展览 ("Exhibitions"), Sheet 4 = 全部类型 ("All types") — the
# latter mirrors the former's rows (interleaved with 演出 rows), so the
# same event/value pairs are applied to both sheets, at their respective
# row numbers.

$wb = $excel.ActiveWorkbook

# Row => new value for column F on the "展览" sheet (Worksheets index 1)
$sheet1Updates = @{
    2  = 220
    3  = 1068
    5  = 389
    6  = 71
    8  = 45
    9  = 6643
    10 = 139
    15 = 1062
    16 = 15911
    19 = 313
    20 = 164
    21 = 112
    22 = 11216
    23 = 820
    24 = 4407
    25 = 282
    27 = 35
    28 = 315
    29 = 134
}

# Row => new value for column F on the "全部类型" sheet (Worksheets index 4)
$sheet4Updates = @{
    2  = 220
    4  = 1068
    6  = 389
    7  = 71
    10 = 45
    11 = 6643
    12 = 139
    18 = 1062
    19 = 15911
    22 = 313
    23 = 164
    24 = 112
    26 = 11216
    27 = 820
    28 = 4407
    29 = 282
    31 = 35
    32 = 315
    33 = 134
}

$ws1 = $wb.Worksheets.Item(1)
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item(4)
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
